# Apply the two changes captured by the commit:
#   1. The table on slide 16 switches from the deck's custom "Table_0"
#      style to the built-in table style {255F9AB2-A36E-4761-8E0C-A2BD64FBFEA5}.
#   2. The presentation's theme colour scheme (ppt/theme/theme1.xml, shared
#      by every slide through the slide master) is swapped from the old
#      "Integral" palette to the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{255F9AB2-A36E-4761-8E0C-A2BD64FBFEA5}")

# --- 2. Theme colours -------------------------------------------------------
# Theme colour scheme index order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# Values below are the "Office Theme" srgbClr values re-packed as COM BGR ints.
$officeThemeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501   # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i]
}
